$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$target = "System, dnasr281@gmail.com"
$replacement = "dnasr281@gmail.com, System"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $target) {
        $cell.Value2 = $replacement
    }
}
